# Add time tracking to game and guess submissions; create GuessTime model
# and update game time logic -- edits the "Functionality for any bugs"
# checklist textbox on slide 2 of the deck.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sp = $s.Shapes.Item(1)
$tr = $sp.TextFrame.TextRange

# 1) "Program tests" -> "Program tests (locally)"
$tr.Paragraphs(2).Runs(1).Text = "Program tests (locally)"

# 2) New sub-bullet "End to End test" right after "Program tests (locally)"
#    (leading carriage-return splits a fresh paragraph off the anchor,
#    inheriting its level-1 bullet formatting).
$tr.Paragraphs(2).InsertAfter([char]13 + "End to End test") | Out-Null

# 3) Two new sub-bullets after "DEE Student tests":
#    "Change everyday " and "Same word for everyone"
$tr.Paragraphs(4).InsertAfter([char]13 + "Change everyday " + [char]13 + "Same word for everyone") | Out-Null

# 4) New top-level bullet "Get Domain (am I get sued for,cloudflare" inserted
#    right before "Award mechanism", written as two runs.
$tr.Paragraphs(8).InsertBefore("Get Domain (am I get sued " + [char]13) | Out-Null
$tr.Paragraphs(8).Runs(1).InsertAfter("for,cloudflare") | Out-Null

# 5) Two new sub-bullets "Everything" and "Stats Tab" after "Appearance of
#    the app".
$tr.Paragraphs(12).InsertAfter([char]13 + "Everything" + [char]13 + "Stats Tab") | Out-Null

# Promote the two freshly-added paragraphs to the same outline level as
# their sibling sub-bullets (must happen last -- re-touching run text
# after this resets the paragraph's level).
$tr.Paragraphs(13).IndentLevel = 2
$tr.Paragraphs(14).IndentLevel = 2
